$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B,D,E,F,G,H,J,L,N across rows 2-25
# (Case with 380 kV line results updated)
$data = @{
    2 = @{ "B"=1.766827342824513; "D"=0.1444903139226028; "E"=0.1701876787727024; "F"=2.082208593905406; "G"=1.567686896685643; "H"=1.338365898435342; "J"=0.2299878191940579; "L"=0.6392872009984103; "N"=1.640143406601247 }
    3 = @{ "B"=1.695111190300565; "D"=0.1397123176989297; "E"=0.1650353036674908; "F"=2.065596191870782; "G"=1.531567113367856; "H"=1.328741585860485; "J"=0.2233584422936801; "L"=0.5877247817050204; "N"=1.652065882509923 }
    4 = @{ "B"=1.651791193987236; "D"=0.1367430120217676; "E"=0.1618467023701378; "F"=2.056828950361734; "G"=1.510619120153137; "H"=1.323692309248685; "J"=0.2192679301154072; "L"=0.5562288416835202; "N"=1.660073885600184 }
    5 = @{ "B"=1.634318017731033; "D"=0.1355239508610708; "E"=0.160541026556686; "F"=2.053615248963609; "G"=1.502390058170363; "H"=1.321850296835123; "J"=0.2175960325786619; "L"=0.5434352689427442; "N"=1.663510372347339 }
    6 = @{ "B"=1.631427502101843; "D"=0.1353209776630777; "E"=0.1603238398253986; "F"=2.05310326312771; "G"=1.501042143409734; "H"=1.321557434953519; "J"=0.2173181162548303; "L"=0.5413134034861287; "N"=1.664091464915074 }
    7 = @{ "B"=1.651554814811078; "D"=0.1367266080808918; "E"=0.161829119058261; "F"=2.056784157303412; "G"=1.510506897818402; "H"=1.323666595129851; "J"=0.2192454024041126; "L"=0.5560561357563643; "N"=1.660119529788055 }
    8 = @{ "B"=1.741951708421936; "D"=0.1428501910162367; "E"=0.1684163298037973; "F"=2.076182467397317; "G"=1.554976467052995; "H"=1.334868512742958; "J"=0.2277062051765881; "L"=0.6214745903312178; "N"=1.644111732575816 }
    9 = @{ "B"=1.924873167651413; "D"=0.1545807988491674; "E"=0.1811360420685482; "F"=2.125654535519573; "G"=1.652026650539312; "H"=1.363693308508914; "J"=0.2441368870358502; "L"=0.7510591552360211; "N"=1.618164783355766 }
    10 = @{ "B"=2.062712828582676; "D"=0.16303721574684; "E"=0.1903625778310172; "F"=2.169066132790476; "G"=1.729464060088446; "H"=1.389102720874234; "J"=0.2561091069970018; "L"=0.8470693715886739; "N"=1.60240656707137 }
    11 = @{ "B"=2.126169686947662; "D"=0.1668506844434887; "E"=0.1945346907089274; "F"=2.19037042422228; "G"=1.766053877020909; "H"=1.401592494686298; "J"=0.2615338090861599; "L"=0.8909250821590717; "N"=1.595952691877443 }
    12 = @{ "B"=2.150307204438491; "D"=0.1682900600958277; "E"=0.1961109749313934; "F"=2.198663185837262; "G"=1.780107729020756; "H"=1.406456772800993; "J"=0.2635848669658714; "L"=0.9075580807889025; "N"=1.593611335277274 }
    13 = @{ "B"=2.145103972051913; "D"=0.1679802723134287; "E"=0.1957716541727805; "F"=2.196867147082571; "G"=1.777072137672292; "H"=1.405403160895133; "J"=0.2631432764682131; "L"=0.9039747215275895; "N"=1.594111028449646 }
    14 = @{ "B"=2.128153337635069; "D"=0.1669691964032012; "E"=0.1946644447781836; "F"=2.19104815082477; "G"=1.767206114118437; "H"=1.401989978074624; "J"=0.261702614529554; "L"=0.8922929726686561; "N"=1.595758012052315 }
    15 = @{ "B"=2.117784606970474; "D"=0.166349273854081; "E"=0.1939857777565166; "F"=2.18751323223114; "G"=1.761188747610959; "H"=1.399916867501844; "J"=0.2608197541980815; "L"=0.885140919903364; "N"=1.596780191759706 }
    16 = @{ "B"=2.05858088224312; "D"=0.1627873331012069; "E"=0.1900894158902418; "F"=2.167705290048318; "G"=1.72710043541673; "H"=1.388305285834122; "J"=0.255754151373182; "L"=0.8442069132049141; "N"=1.602842707685895 }
    17 = @{ "B"=2.022453801932102; "D"=0.1605937148530927; "E"=0.187692700986922; "F"=2.155953397443028; "G"=1.70653895637264; "H"=1.381420955318077; "J"=0.25264100994427; "L"=0.8191412878015853; "N"=1.606744767191174 }
    18 = @{ "B"=2.001745355814649; "D"=0.1593288540152713; "E"=0.1863118126654228; "F"=2.149340398835108; "G"=1.694840718989013; "H"=1.377548842573646; "J"=0.2508483914714219; "L"=0.8047411688459647; "N"=1.609056401491074 }
    19 = @{ "B"=1.994746023077596; "D"=0.1589000498804154; "E"=0.1858438619849068; "F"=2.147126445940145; "G"=1.690901849613851; "H"=1.376252827855865; "J"=0.2502410963830641; "L"=0.7998684536503617; "N"=1.609850640831084 }
    20 = @{ "B"=2.02629225750934; "D"=0.1608275545129914; "E"=0.1879480797485726; "F"=2.1571892458512; "G"=1.708714479747783; "H"=1.382144734857746; "J"=0.2529726188888191; "L"=0.8218078114031755; "N"=1.606322425113234 }
    21 = @{ "B"=2.133129228446819; "D"=0.1672663007171025; "E"=0.1949897568116299; "F"=2.192751204611568; "G"=1.770098611937442; "H"=1.402988851143647; "J"=0.2621258581073818; "L"=0.895723486343087; "N"=1.595271470235105 }
    22 = @{ "B"=2.203581390380918; "D"=0.17144703322775; "E"=0.1995708927592972; "F"=2.217306874357746; "G"=1.8113722405499; "H"=1.417397040645426; "J"=0.268089586512005; "L"=0.9441819472651787; "N"=1.588646907731786 }
    23 = @{ "B"=2.165922437040251; "D"=0.1692181675601176; "E"=0.1971277752704026; "F"=2.204080325076433; "G"=1.789237307706202; "H"=1.409634988314735; "J"=0.2649083409073398; "L"=0.9183050383549016; "N"=1.592127910813758 }
    24 = @{ "B"=2.024556700166443; "D"=0.1607218471990137; "E"=0.1878326323082504; "F"=2.156630072487573; "G"=1.707730543194913; "H"=1.381817247002004; "J"=0.2528227073010214; "L"=0.8206022435457214; "N"=1.60651315325066 }
    25 = @{ "B"=1.874782806369581; "D"=0.1514363459239405; "E"=0.1777160001470435; "F"=2.111037378091936; "G"=1.624704350993881; "H"=1.355156015436023; "J"=0.2397093183505206; "L"=0.7158627679522169; "N"=1.624602758834413 }
}

foreach ($row in $data.Keys) {
    $rowVals = $data[$row]
    foreach ($col in $rowVals.Keys) {
        $ws.Range("$col$row").Value2 = $rowVals[$col]
    }
}
